# Update the correlation matrix and the three "Toggles" count sheets to
# reflect results from a gadget run that now works "with any number of
# inputs, with full simulation or partial" (see commit message). The
# three toggle-count sheets feed the bar charts on each of those tabs.

$wb = $excel.ActiveWorkbook

# --- "Correlation matrix": row "no delays" goes to 0 correlation on all
# inputs, while "gate delays" / "gate+inputs delay" both settle on the
# same correlation value across every input column.
$wsCorr = $wb.Worksheets.Item("Correlation matrix")
$wsCorr.Range("B2:E2").Value = 0
$wsCorr.Range("B3:E3").Value = 0.1796053020267749
$wsCorr.Range("B4:E4").Value = 0.1796053020267749

# --- "Toggles no del": no-delay toggle counts for the two HW/HD inputs.
$wsNoDel = $wb.Worksheets.Item("Toggles no del")
$wsNoDel.Range("B2").Value = 128
$wsNoDel.Range("B3").Value = 128

# --- "Toggles del": gate-delay toggle counts.
$wsDel = $wb.Worksheets.Item("Toggles del")
$wsDel.Range("B2").Value = 80
$wsDel.Range("B3").Value = 128
$wsDel.Range("B4").Value = 48

# --- "Toggles input del": gate+input-delay toggle counts.
$wsInputDel = $wb.Worksheets.Item("Toggles input del")
$wsInputDel.Range("B2").Value = 80
$wsInputDel.Range("B3").Value = 128
$wsInputDel.Range("B4").Value = 48
